function Set-TextValue {
    param($Sheet, $CellRef, $NewValue)
    $rng = $Sheet.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $NewValue
    $rng.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws "D2" "67.934.86"
Set-TextValue $ws "E2" "  -0.23%  "
Set-TextValue $ws "D3" "3.804.47"
Set-TextValue $ws "E3" "  -1.68%  "
Set-TextValue $ws "E4" "  -0.03%  "
Set-TextValue $ws "D5" "599.31"
Set-TextValue $ws "E5" "  -0.02%  "
Set-TextValue $ws "D6" "169.37"
Set-TextValue $ws "E6" "  +1.22%  "
Set-TextValue $ws "D7" "3.809.10"
Set-TextValue $ws "E7" "  -1.54%  "
Set-TextValue $ws "E8" "  -0.03%  "
Set-TextValue $ws "D9" "0.532"
Set-TextValue $ws "E9" "  +0.92%  "
Set-TextValue $ws "E10" "  +0.38%  "
Set-TextValue $ws "D11" "6.49"
Set-TextValue $ws "E11" "  +1.23%  "
Set-TextValue $ws "D12" "0.461"
Set-TextValue $ws "E12" "  +0.72%  "
Set-TextValue $ws "E13" "  +12.00%  "
Set-TextValue $ws "D14" "36.92"
Set-TextValue $ws "E14" "  -0.07%  "
Set-TextValue $ws "D15" "4.438.44"
Set-TextValue $ws "E15" "  -1.80%  "
Set-TextValue $ws "D16" "3.825.27"
Set-TextValue $ws "E16" "  -1.24%  "
Set-TextValue $ws "D17" "67.975.43"
Set-TextValue $ws "E17" "  -0.23%  "
Set-TextValue $ws "D18" "18.26"
Set-TextValue $ws "E18" "  -0.03%  "
Set-TextValue $ws "D19" "7.44"
Set-TextValue $ws "E19" "  +0.35%  "
Set-TextValue $ws "D21" "10.88"
Set-TextValue $ws "D22" "470.73"
Set-TextValue $ws "E22" "  +1.06%  "
Set-TextValue $ws "D23" "0.726"
Set-TextValue $ws "E23" "  -0.73%  "
Set-TextValue $ws "D24" "0.0000151"
Set-TextValue $ws "E24" "  -7.20%  "
Set-TextValue $ws "D25" "83.42"
Set-TextValue $ws "E25" "  -0.09%  "
Set-TextValue $ws "D26" "2.29"
Set-TextValue $ws "E26" "  +1.94%  "
Set-TextValue $ws "D27" "12.24"
Set-TextValue $ws "E27" "  +0.99%  "
Set-TextValue $ws "D28" "10.34"
Set-TextValue $ws "E28" "  +3.31%  "
Set-TextValue $ws "E29" "  +0.03%  "
Set-TextValue $ws "E30" "  -0.62%  "
Set-TextValue $ws "D31" "3.949.21"
Set-TextValue $ws "E31" "  -1.74%  "
Set-TextValue $ws "D32" "7.75"
Set-TextValue $ws "E32" "  -1.29%  "
Set-TextValue $ws "D33" "2.29"
Set-TextValue $ws "E33" "  -0.69%  "
Set-TextValue $ws "D34" "30.88"
Set-TextValue $ws "E34" "  -1.15%  "
Set-TextValue $ws "D35" "9.36"
Set-TextValue $ws "E35" "  +0.09%  "
Set-TextValue $ws "D36" "3.767.58"
Set-TextValue $ws "E36" "  -1.96%  "
Set-TextValue $ws "B37" "dogwifhat"
Set-TextValue $ws "C37" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws "D37" "3.86"
Set-TextValue $ws "E37" "  +12.55%  "
Set-TextValue $ws "B38" "Hedera"
Set-TextValue $ws "C38" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws "D38" "0.107"
Set-TextValue $ws "E38" "  +2.95%  "
Set-TextValue $ws "D39" "0.140"
Set-TextValue $ws "E39" "  +0.12%  "
Set-TextValue $ws "D40" "1.02"
Set-TextValue $ws "E40" "  -0.87%  "
Set-TextValue $ws "D41" "5.95"
Set-TextValue $ws "E41" "  +0.74%  "
Set-TextValue $ws "D42" "0.999"
Set-TextValue $ws "E42" "  -0.18%  "
Set-TextValue $ws "D43" "0.318"
Set-TextValue $ws "E43" "  +1.73%  "
Set-TextValue $ws "B44" "USDe"
Set-TextValue $ws "C44" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws "D44" "1.00"
Set-TextValue $ws "E44" "  +0.00%  "
Set-TextValue $ws "B45" "Stacks"
Set-TextValue $ws "C45" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws "D45" "1.98"
Set-TextValue $ws "E45" "  -0.18%  "
Set-TextValue $ws "B46" "Cosmos"
Set-TextValue $ws "C46" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws "D46" "8.76"
Set-TextValue $ws "E46" "  +2.71%  "
Set-TextValue $ws "D47" "0.000290"
Set-TextValue $ws "E47" "  +4.19%  "
Set-TextValue $ws "B48" "Bittensor"
Set-TextValue $ws "C48" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws "D48" "409.40"
Set-TextValue $ws "E48" "  -4.58%  "
Set-TextValue $ws "B49" "OKB"
Set-TextValue $ws "C49" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws "D49" "46.49"
Set-TextValue $ws "E49" "  -1.83%  "
Set-TextValue $ws "D50" "141.85"
Set-TextValue $ws "E50" "  -1.65%  "
Set-TextValue $ws "D51" "0.0359"
Set-TextValue $ws "E51" "  +0.61%  "

$wb.Save()
